$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10/11: Cardano and Dogecoin swap places (with updated price/volume).
# D values are numeric-looking so a leading apostrophe forces text storage,
# matching the original inlineStr cell type instead of being coerced to a number.
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "'0.164"
$ws.Range("E10").Value = "  +8.62%  "
$ws.Range("B11").Value = "Cardano"
$ws.Range("C11").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D11").Value = "'0.630"
$ws.Range("E11").Value = "  +2.37%  "

# Remaining price / volume updates
$ws.Range("D2").Value = "64.492.76"
$ws.Range("E2").Value = "  +0.95%  "
$ws.Range("D3").Value = "3.353.86"
$ws.Range("E3").Value = "  -1.34%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'560.00"
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("D6").Value = "'175.95"
$ws.Range("E6").Value = "  +2.47%  "
$ws.Range("E7").Value = "  +1.17%  "
$ws.Range("D8").Value = "3.345.55"
$ws.Range("E8").Value = "  -1.50%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D12").Value = "'54.97"
$ws.Range("E12").Value = "  -1.63%  "
$ws.Range("E13").Value = "  +3.46%  "
$ws.Range("D14").Value = "'9.08"
$ws.Range("E14").Value = "  +1.59%  "
$ws.Range("D15").Value = "3.899.25"
$ws.Range("E15").Value = "  -1.24%  "
$ws.Range("D16").Value = "'18.26"
$ws.Range("E16").Value = "  +2.01%  "
$ws.Range("D17").Value = "3.355.98"
$ws.Range("E17").Value = "  -4.14%  "
$ws.Range("E18").Value = "  -1.64%  "
$ws.Range("D19").Value = "'11.82"
$ws.Range("E19").Value = "  +0.99%  "
$ws.Range("D20").Value = "64.474.52"
$ws.Range("E20").Value = "  +0.60%  "
$ws.Range("D21").Value = "'0.987"
$ws.Range("E21").Value = "  +0.38%  "
$ws.Range("D22").Value = "'462.10"
$ws.Range("E22").Value = "  +13.48%  "
$ws.Range("D23").Value = "'4.85"
$ws.Range("E23").Value = "  +12.33%  "
$ws.Range("D24").Value = "'4.10"
$ws.Range("E24").Value = "  +0.35%  "
$ws.Range("D25").Value = "'86.22"
$ws.Range("E25").Value = "  +4.32%  "
$ws.Range("D26").Value = "'13.44"
$ws.Range("E26").Value = "  +1.70%  "
$ws.Range("D27").Value = "'10.84"
$ws.Range("E27").Value = "  +1.07%  "
$ws.Range("D28").Value = "'2.84"
$ws.Range("E28").Value = "  +3.27%  "
$ws.Range("D29").Value = "'8.78"
$ws.Range("E29").Value = "  +0.91%  "
$ws.Range("D30").Value = "'30.16"
$ws.Range("E30").Value = "  +2.47%  "
$ws.Range("D31").Value = "'6.67"
$ws.Range("E31").Value = "  +0.80%  "
$ws.Range("D32").Value = "'11.48"
$ws.Range("E32").Value = "  +0.55%  "
$ws.Range("D33").Value = "'580.25"
$ws.Range("E33").Value = "  -1.19%  "
$ws.Range("E34").Value = "  +1.32%  "
$ws.Range("D35").Value = "'59.17"
$ws.Range("E35").Value = "  +0.96%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("E37").Value = "  -6.28%  "
$ws.Range("D38").Value = "'35.79"
$ws.Range("E38").Value = "  -0.20%  "
$ws.Range("D39").Value = "'3.48"
$ws.Range("E39").Value = "  +2.65%  "
$ws.Range("D40").Value = "0.0₃0755"
$ws.Range("E40").Value = "  +4.23%  "
$ws.Range("D41").Value = "'0.371"
$ws.Range("E41").Value = "  +0.61%  "
$ws.Range("D42").Value = "3.089.87"
$ws.Range("E42").Value = "  -2.86%  "
$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = "  -0.37%  "
$ws.Range("E44").Value = "  -1.22%  "
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").Value = "'0.0412"
$ws.Range("E46").Value = "  +2.25%  "
$ws.Range("E47").Value = "  -0.59%  "
$ws.Range("E48").Value = "  +2.16%  "
$ws.Range("E49").Value = "  -0.48%  "
$ws.Range("D50").Value = "'8.34"
$ws.Range("E50").Value = "  +1.62%  "
$ws.Range("D51").Value = "'136.09"
$ws.Range("E51").Value = "  +0.82%  "
